# Add a new "Turtle Count" column, inserted before the existing
# "Turtle Activity" column (U), which shifts "Turtle Activity" to column V.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at U; this pushes the old column U ("Turtle
# Activity") one position right to V, copying formatting/styles and
# updating the sheet dimension + row spans automatically.
$ws.Columns("U").Insert() | Out-Null

# New column header + data values ("Turtle Count").
$ws.Range("U1").Value2 = "Turtle Count"

$turtleCounts = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 14
    6  = 13
    7  = 12
    8  = 1
    9  = 8
    10 = 8
    11 = 7
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 3
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 2
    27 = 3
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 0
}

foreach ($row in $turtleCounts.Keys) {
    $ws.Cells.Item($row, 21).Value2 = $turtleCounts[$row]
}

# Re-apply the AutoFilter over the now-wider used range (A1:V64) -- the
# column insert alone does not grow the existing filter range.
$ws.AutoFilterMode = $false
$ws.Range("A1:V64").AutoFilter() | Out-Null

# Un-hide the _FilterDatabase defined name and point it at the new range
# (Excel clears the hidden flag once the filter is (re)applied through the UI).
$fd = $ws.Names.Item(1)
$fd.RefersTo = "=Sheet1!`$A`$1:`$V`$64"
$fd.Visible = $true

# Match the selection left behind by selecting the whole new column U.
$ws.Columns("U").Select() | Out-Null
